$d = $word.ActiveDocument

# ============================================================
# Part 1: Insert a new "Meta description" paragraph right after
# the first (title) paragraph.
# ============================================================
$p1 = $d.Paragraphs.Item(1)
$insertPoint = $d.Range($p1.Range.End, $p1.Range.End)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t xml:space="preserve">: Check out our review of 4 Secret Pyramids free play slot game, with info on bonus features, ways to win, and high volatility.</w:t></w:r></w:p>
          <w:p/>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertPoint.InsertXML($xml) | Out-Null

# The InsertXML call above inserted two paragraph marks (one that
# ends our new paragraph, and an extra empty one). Remove the extra
# paragraph mark so the following paragraph ("Graphics and User
# Experience") merges back to being immediately after our new
# paragraph, keeping its own original formatting intact.
$p2 = $d.Paragraphs.Item(2)
$extraMark = $d.Range($p2.Range.End - 1, $p2.Range.End)
$extraMark.Delete()

# ============================================================
# Part 2: Remove the duplicated bold heading paragraph further
# down ("Play 4 Secret Pyramids Free - Review of Bonuses and
# Features") that precedes the italic meta-description paragraph.
# Search begins after our newly-inserted paragraph (p2) so the
# occurrence of this title text at the very top of the document
# is not matched instead.
# ============================================================
$p2 = $d.Paragraphs.Item(2)
$afterNewPara = $d.Range($p2.Range.End, $d.Content.End)
$afterNewPara.Find.Execute("Play 4 Secret Pyramids Free - Review of Bonuses and Features") | Out-Null
$dupHeading = $d.Range($afterNewPara.Start, $afterNewPara.End)
$dupHeading.Expand(4) | Out-Null   # wdParagraph
$dupHeading.Delete()

# ============================================================
# Part 3: Replace the text of the remaining italic paragraph with
# the new image-prompt text, preserving its italic formatting.
# Search begins right after our inserted meta-description
# paragraph, so the earlier (new) paragraph is not matched.
# ============================================================
$italicSearch = $d.Range($p2.Range.End, $d.Content.End)
$italicSearch.Find.Execute("Check out our review of 4 Secret Pyramids free play slot game, with info on bonus features, ways to win, and high volatility.") | Out-Null
$targetRange = $d.Range($italicSearch.Start, $italicSearch.End)
$targetRange.Text = "Please create a cartoon-style feature image for the 4 Secret Pyramids game. The image should prominently feature a happy Maya warrior wearing glasses. The cartoon character should be holding up a pyramid with a big smile on their face, with other pyramids visible in the background. The image should be bright and colorful, conveying the excitement and fun of the game. It should be eye-catching and make viewers want to play the game. The image should be in a landscape orientation, suitable for use on a website or in social media posts."
